# Update column F ("dSF") values on Sheet1 per the source-data repull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = -5
    4  = -2
    5  = -3
    6  = 4
    9  = 4
    10 = 1
    11 = 1
    12 = 2
    13 = 8
    14 = 4
    15 = -2
    16 = -3
    17 = -3
    18 = -3
    19 = 2
    20 = -1
    21 = -1
    22 = -2
    23 = 2
    24 = 1
    25 = -3
    26 = 5
    27 = -7
    28 = -3
    30 = -1
    31 = 1
    32 = -3
    33 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
